$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "50.039.15"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +4.18%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.660.20"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +7.02%  "

$ws.Range("E4").Value = "  +0.09%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "114.17"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +8.25%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "326.50"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.90%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.530"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.01%  "

$ws.Range("E8").Value = "  +0.04%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.557"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +3.81%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "41.35"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +6.39%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.14"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.04%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0826"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +3.18%  "

$ws.Range("E13").Value = "  +0.46%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.39"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +4.15%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.077.50"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +7.01%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.648.40"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +6.41%  "

$ws.Range("E17").Value = "  +6.15%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "49.994.48"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +4.29%  "

$ws.Range("E19").Value = "  +4.63%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.80"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.67%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.95"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.55%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.0₃0962"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.49%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "72.54"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.21%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "276.38"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.64%  "

$ws.Range("E25").Value = "  +3.74%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.96"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +4.96%  "

$ws.Range("E27").Value = "  +0.09%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.04"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.44%  "

$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.23"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.56%  "

$ws.Range("B30").Value = "InjectiveProtocol"
$ws.Range("C30").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "36.79"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +6.42%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.142"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.56%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "50.25"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.88%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.51"
$ws.Range("D33").Style = "Normal"

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "19.70"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.57%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0816"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +5.83%  "

$ws.Range("E36").Value = "  -0.09%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.04"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +10.53%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.08"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +7.84%  "

$ws.Range("E39").Value = "  +9.64%  "

$ws.Range("B40").Value = "Stellar"
$ws.Range("C40").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.113"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.72%  "

$ws.Range("B41").Value = "Monero"
$ws.Range("C41").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "124.37"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.37%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "22.43"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.54%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.22"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.21%  "

$ws.Range("E44").Value = "  +5.24%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.108.11"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +5.48%  "

$ws.Range("E46").Value = "  +5.63%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.26"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +13.36%  "

$ws.Range("E48").Value = "  +5.33%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.12"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.43%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "5.37"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.75%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "59.97"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +6.11%  "
